$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Professor names (column A) ---
$ws.Range("A2").Value = "Prof. Smith"
$ws.Range("A3").Value = "Prof. Cole"

# --- Professor e-mails (column B) -- Prof. Cole's becomes a mailto hyperlink ---
$ws.Range("B3").Value = "John.Cole@utdallas.edu "
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:John.Cole@utdallas.edu", "", "", "John.Cole@utdallas.edu ")
$ws.Range("B2").Value = "Jason.Smith1@utdallas.edu"

# --- Section (column D), stored as a literal quoted string "001" ---
$ws.Range("D2").Value = '"001"'
$ws.Range("D3").Value = '"001"'

# --- Course name (column E) ---
$ws.Range("E2").Value = "Computer Science II"
$ws.Range("E3").Value = "Computer Architecture"

# --- Recommended student name/netid (columns F, G) ---
$ws.Range("F2").Value = "John Doe"
$ws.Range("G2").Value = "jde200000"
$ws.Range("F3").Value = "Ron Doe"
$ws.Range("G3").Value = "rde200001"

# --- Keywords (column I) ---
$ws.Range("I2").Value = "C++"
$ws.Range("I3").Value = "C"

# --- Course number (column C) ---
$ws.Range("C2").Value = 2337
$ws.Range("C3").Value = 2340
$ws.Range("C3").Font.Color = 1381653

# --- Num of graders (column H) ---
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 2

# --- Column widths, widened to fit the new, longer content ---
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Final selection, as left by the editing session ---
$ws.Range("I4").Select()
